$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (2-530).
# The diff shows every one of these cells moving from serial 45192
# (2023-09-23) to serial 45202 (2023-10-03).
$ws.Range("C2:C530").Value = 45202
